# Insert two new data rows right after the current row 175 (i.e. at rows 176-177),
# pushing the existing rows 176.. down by two (they become 178.. automatically).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(176).Insert()
$ws.Rows(176).Insert()

# New row 176
$ws.Cells.Item(176, 1).Value = 6
$ws.Cells.Item(176, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(176, 3).Value = "Metropolitana"
$ws.Cells.Item(176, 4).Value = 44460
$ws.Cells.Item(176, 5).Value = 13
$ws.Cells.Item(176, 6).Value = 100112043
$ws.Cells.Item(176, 7).Value = "Pepino ensalada"
$ws.Cells.Item(176, 8).Value = "Sin especificar"
$ws.Cells.Item(176, 9).Value = "Primera"
$ws.Cells.Item(176, 10).Value = 400
$ws.Cells.Item(176, 11).Value = 14000
$ws.Cells.Item(176, 12).Value = 15000
$ws.Cells.Item(176, 13).Value = 14425
$ws.Cells.Item(176, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(176, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(176, 16).Value = 240
$ws.Cells.Item(176, 17).Value = 60
$ws.Cells.Item(176, 18).Value = "Hortaliza"

# New row 177
$ws.Cells.Item(177, 1).Value = 6
$ws.Cells.Item(177, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(177, 3).Value = "Metropolitana"
$ws.Cells.Item(177, 4).Value = 44460
$ws.Cells.Item(177, 5).Value = 13
$ws.Cells.Item(177, 6).Value = 100112043
$ws.Cells.Item(177, 7).Value = "Pepino ensalada"
$ws.Cells.Item(177, 8).Value = "Sin especificar"
$ws.Cells.Item(177, 9).Value = "Segunda"
$ws.Cells.Item(177, 10).Value = 220
$ws.Cells.Item(177, 11).Value = 12000
$ws.Cells.Item(177, 12).Value = 13000
$ws.Cells.Item(177, 13).Value = 12545
$ws.Cells.Item(177, 14).Value = "$/caja 100 unidades"
$ws.Cells.Item(177, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(177, 16).Value = 125
$ws.Cells.Item(177, 17).Value = 100
$ws.Cells.Item(177, 18).Value = "Hortaliza"
